# Login Using DataDriven Completed
# Trim the login-data table down to the two rows actually used by the
# data-driven test (standard_user / locked_out_user), normalize the
# expected-result strings to "Pass"/"Fail", clear the now-unused rows
# (keeping their formatting), set the page to portrait orientation, and
# move the active selection off the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "expected" column for the two rows that remain meaningful.
$ws.Range("C2").Value = "Pass"
$ws.Range("C3").Value = "Fail"

# Rows 4-7 no longer carry user data; clear their contents but keep the
# existing cell formatting/style in place.
$ws.Range("A4:C7").ClearContents()

# Flip the sheet to portrait orientation.
$ws.PageSetup.Orientation = 1

# Move the selection away from the data (matches the saved selection in
# the workbook after the edit).
$ws.Range("D1:F7").Select() | Out-Null
